$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 1, 5, 5),
    @(8, 1, 10, 10),
    @(10, 1, 15, 15),
    @(6, 2, 5, 5),
    @(5, 3, 5, 5),
    @(2, 4, 5, 5),
    @(3, 4, 10, 10),
    @(11, 4, 15, 15),
    @(12, 4, 20, 21),
    @(1, 5, 5, 5),
    @(7, 5, 10, 10),
    @(9, 5, 15, 16),
    @(10, 5, 21, 21)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
